# Add 6 newly-reported COVID deaths (dated 43941 / 2020-04-20) plus fill in
# the already-numbered-but-previously-blank rows 172-173 on the "10yr" sheet.
# (The record that used to live in row 174 - age 87, woman, Me'ainei
# HaYeshua - is now the first of this batch, in row 172.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10yr")

$green = 6137344   # RGB(0,166,93) - the "new entry today" highlight color
$dateFmt = "MM/DD/YY"

function Set-DeathRow {
    param($row, $num, $date, $age, $ageIsNew, $comment, $sex, $place)

    $ws.Range("B$row").Value = $num

    $ws.Range("A$row").Value = $date
    $ws.Range("A$row").NumberFormat = $dateFmt

    $ws.Range("C$row").Value = $age
    if ($ageIsNew) {
        $ws.Range("C$row").Font.Color = $green
    }

    if ($comment -ne "") {
        $ws.Range("D$row").Value = $comment
    }

    $ws.Range("E$row").Value = $sex
    $ws.Range("F$row").Value = $place
}

# row172: pre-existing record (age 87, w, Me'ainei HaYeshua) moved down into
# the slot that already had B172=171
Set-DeathRow 172 171 43940 87 $false "" "w" "מעייני הישועה"

# row173: new record
Set-DeathRow 173 172 43941 90 $true "" "w" "מרכז רפואי שהם"

# row174: new record
Set-DeathRow 174 173 43941 89 $true "מחלות רקע רבות" "m" "שמיר אסף הרופא"

# row175: new record
Set-DeathRow 175 174 43941 62 $false "" "m" "שיבא"

# row176: new record
Set-DeathRow 176 175 43941 71 $true "מחלות רקע שונות" "m" "הלל יפה"

# row177: new record
Set-DeathRow 177 176 43941 63 $true "מחלות רקע" "m" "הדסה עין כרם"

# row178: new record
Set-DeathRow 178 177 43941 85 $true "מחלות רקע" "w" "הדסה עין כרם"

# Update the worksheet view to scroll to the newly-added rows, like the
# author did before saving.
$ws.Application.ActiveWindow.ScrollRow = 167
$ws.Range("C173").Select()
